$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.391.55'
$ws.Range("E2").Value = '  +2.24%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.099.93'
$ws.Range("E3").Value = '  +0.21%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.83%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '343.55'
$ws.Range("E5").Value = '  -0.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5357'
$ws.Range("E7").Value = '  +3.74%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4441'
$ws.Range("E8").Value = '  +1.31%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '54.83'
$ws.Range("E9").Value = '  +3.47%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09412'
$ws.Range("E10").Value = '  +1.83%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.172'
$ws.Range("E11").Value = '  +0.59%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.78'
$ws.Range("E12").Value = '  -0.46%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.160.57'
$ws.Range("E13").Value = '  +2.91%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.567'
$ws.Range("E14").Value = '  +4.04%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.919'
$ws.Range("E15").Value = '  +2.16%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '102.01'
$ws.Range("E16").Value = '  +2.43%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001162'
$ws.Range("E17").Value = '  +0.84%  '

$ws.Range("E18").Value = '  -0.73%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '21.20'
$ws.Range("E19").Value = '  +1.83%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.06691'
$ws.Range("E20").Value = '  +0.48%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.331'
$ws.Range("E21").Value = '  +2.02%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  -0.70%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.400.82'
$ws.Range("E23").Value = '  +2.15%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.55'
$ws.Range("E24").Value = '  +0.30%  '

$ws.Range("E25").Value = '  +0.09%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '21.92'
$ws.Range("E26").Value = '  -0.13%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.90'
$ws.Range("E27").Value = '  +0.95%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.527'
$ws.Range("E28").Value = '  +0.27%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.791'
$ws.Range("E29").Value = '  +7.68%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.74'
$ws.Range("E30").Value = '  +0.48%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.145'
$ws.Range("E31").Value = '  +0.44%  '

$ws.Range("E32").Value = '  +0.71%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.668'
$ws.Range("E33").Value = '  +1.76%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.265'
$ws.Range("E34").Value = '  +1.46%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.848'
$ws.Range("E35").Value = '  -2.64%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.16'
$ws.Range("E36").Value = '  -0.47%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02638'
$ws.Range("E37").Value = '  +2.41%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06805'
$ws.Range("E38").Value = '  +1.15%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.7027'
$ws.Range("E39").Value = '  -0.98%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.68'
$ws.Range("E40").Value = '  +1.65%  '

$ws.Range("E41").Value = '  +2.06%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2223'
$ws.Range("E42").Value = '  -0.29%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6869'
$ws.Range("E43").Value = '  -1.70%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.46'
$ws.Range("E44").Value = '  +1.13%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.343'
$ws.Range("E45").Value = '  +1.03%  '

$ws.Range("E46").Value = '  -0.64%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.398'
$ws.Range("E47").Value = '  +20.39%  '

$ws.Range("E48").Value = '  +0.52%  '

$ws.Range("E49").Value = '  +9.41%  '

$ws.Range("E50").Value = '  -2.53%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.221'
$ws.Range("E51").Value = '  +0.08%  '

